# Update "want to go" counts (column F) on both the "展览" and "全部类型"
# worksheets, which carry duplicate data in this workbook.

$wb = $excel.ActiveWorkbook

# Row number (key) -> new F-column value (value)
$updates = @{
    2  = 10
    3  = 1094
    5  = 87
    7  = 58
    8  = 11260
    9  = 4292
    11 = 25
    13 = 2506
    15 = 110
    16 = 19
    17 = 165
    18 = 491
    19 = 11251
    20 = 11099
    22 = 37
    25 = 36
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
